$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains values that look numeric (e.g. "0.9989"),
# as well as "thousands"-style values with multiple dots (e.g. "24.677.05").
# In the source workbook these are stored as plain text, so force the
# column to Text format before writing the new values to avoid Excel
# silently converting them to numbers (which would also mangle values
# like "0.00001315" into scientific notation). The format is reset back
# to Normal afterwards so no visible formatting change is introduced.
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = "24.677.05"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "1.695.42"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "316.48"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.3943"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "0.4019"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "1.521"
$ws.Range("E9").Value = "  +3.75%  "
$ws.Range("D10").Value = "0.9991"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "52.41"
$ws.Range("E11").Value = "  -3.53%  "
$ws.Range("D12").Value = "0.08759"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "7.225"
$ws.Range("E13").Value = "  +6.22%  "
$ws.Range("D14").Value = "23.33"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "8.142"
$ws.Range("E15").Value = "  +11.49%  "
$ws.Range("D16").Value = "0.00001315"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "1.694.98"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "99.96"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "0.07057"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "19.75"
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("D21").Value = "7.051"
$ws.Range("E21").Value = "  +6.44%  "
$ws.Range("D22").Value = "0.9990"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "14.23"
$ws.Range("E23").Value = "  +3.00%  "
$ws.Range("D24").Value = "24.683.62"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").Value = "3.146"
$ws.Range("E25").Value = "  +8.90%  "
$ws.Range("D26").Value = "2.341"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "22.75"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").Value = "162.03"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "137.29"
$ws.Range("E29").Value = "  +5.00%  "
$ws.Range("D30").Value = "5.183"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "7.545"
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").Value = "1.883.34"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "1.086"
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("D34").Value = "0.08611"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "7.197"
$ws.Range("E35").Value = "  +7.81%  "
$ws.Range("E36").Value = "  +10.76%  "
$ws.Range("D37").Value = "0.2742"
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("D38").Value = "1.919"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "14.47"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "0.09126"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").Value = "0.02735"
$ws.Range("E41").Value = "  +8.15%  "
$ws.Range("D42").Value = "1.481"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").Value = "0.7649"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").Value = "2.654"
$ws.Range("E44").Value = "  +9.57%  "
$ws.Range("D45").Value = "0.7189"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "15.67"
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("D47").Value = "4.221"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Value = "1.324"
$ws.Range("E50").Value = "  +9.11%  "
$ws.Range("D51").Value = "0.08007"
$ws.Range("E51").Value = "  +2.36%  "

# Restore the original (default/Normal) cell style on the Price column so
# that no extraneous formatting is left behind on the cells themselves.
$dCol.Style = "Normal"
